$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A29").Value = "NOT"
$ws.Range("B29").Value = 100

$ws.Range("A30").Value = "ERROR"

$ws.Range("B30").Select()
